$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D (shifts existing D:K data to F:M)
$ws.Columns("D:E").Insert()

# Copy number formats from the (now-shifted) first two old data columns (F:G) into new D:E
# so the new columns inherit the correct date/number formatting, done per contiguous block
# to avoid touching the blank separator rows (36, 78).
$ws.Range("F7:G35").Copy()
$ws.Range("D7:E35").PasteSpecial(-4122)
$ws.Range("F38:G77").Copy()
$ws.Range("D38:E77").PasteSpecial(-4122)
$ws.Range("F80:G102").Copy()
$ws.Range("D80:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the two new columns (D,E) with the newest two quarters of data
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 10900
$ws.Range("E8").Value = 14200
$ws.Range("D9").Value = 5400
$ws.Range("E9").Value = 6600
$ws.Range("D10").Value = 5500
$ws.Range("E10").Value = 7600
$ws.Range("D12").Value = 2700
$ws.Range("E12").Value = 4400
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = "NA"
$ws.Range("E14").Value = "NA"
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("D17").Value = 21700
$ws.Range("E17").Value = 26200
$ws.Range("D18").Value = -10800
$ws.Range("E18").Value = -12000
$ws.Range("D20").Value = -300
$ws.Range("E20").Value = 0
$ws.Range("D21").Value = -10600
$ws.Range("E21").Value = -11500
$ws.Range("D22").Value = 800
$ws.Range("E22").Value = 1000
$ws.Range("D23").Value = -11900
$ws.Range("E23").Value = -13000
$ws.Range("D24").Value = -15000
$ws.Range("E24").Value = 0
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 3100
$ws.Range("E26").Value = -13000
$ws.Range("D27").Value = 3100
$ws.Range("E27").Value = -13000
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = -15900
$ws.Range("E29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = 300
$ws.Range("E32").Value = 0
$ws.Range("D33").Value = -12900
$ws.Range("E33").Value = -13000
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = -12900
$ws.Range("E35").Value = -13000
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 99200
$ws.Range("E41").Value = 101900
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("D43").Value = 12300
$ws.Range("E43").Value = 11400
$ws.Range("D44").Value = 6600
$ws.Range("E44").Value = 4500
$ws.Range("D45").Value = 1100
$ws.Range("E45").Value = 1400
$ws.Range("D46").Value = 119300
$ws.Range("E46").Value = 119100
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 0
$ws.Range("D48").Value = 5200
$ws.Range("E48").Value = 5700
$ws.Range("D49").Value = 33500
$ws.Range("E49").Value = 39400
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 0
$ws.Range("E52").Value = 0
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 158000
$ws.Range("E54").Value = 164200
$ws.Range("D57").Value = 8000
$ws.Range("E57").Value = 200
$ws.Range("D58").Value = 0
$ws.Range("E58").Value = 0
$ws.Range("D59").Value = 12600
$ws.Range("E59").Value = 15400
$ws.Range("D60").Value = 20500
$ws.Range("E60").Value = 15600
$ws.Range("D61").Value = 44100
$ws.Range("E61").Value = 43800
$ws.Range("D62").Value = 500
$ws.Range("E62").Value = 500
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 65100
$ws.Range("E66").Value = 60000
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = -126000
$ws.Range("E72").Value = -113100
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 92900
$ws.Range("E76").Value = 104200
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = -12900
$ws.Range("E81").Value = -13000
$ws.Range("D83").Value = 500
$ws.Range("E83").Value = 500
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = -7300
$ws.Range("E89").Value = -11600
$ws.Range("D91").Value = -200
$ws.Range("E91").Value = -100
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = 4400
$ws.Range("E94").Value = -100
$ws.Range("D96").Value = 0
$ws.Range("E96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = 200
$ws.Range("E100").Value = 69900
$ws.Range("D101").Value = 0
$ws.Range("E101").Value = 0
$ws.Range("D102").Value = -2600
$ws.Range("E102").Value = 58200

# Apply revised historical values in column H (previously column F before the insert)
$ws.Range("H8").Value = 7200
$ws.Range("H9").Value = 4300
$ws.Range("H10").Value = 2900
$ws.Range("H12").Value = 1400
$ws.Range("H17").Value = 15100
$ws.Range("H18").Value = -7900
$ws.Range("H21").Value = -7600
$ws.Range("H23").Value = -8800
$ws.Range("H24").Value = -2800
$ws.Range("H26").Value = -6000
$ws.Range("H27").Value = -6000
$ws.Range("H29").Value = -2600
$ws.Range("H43").Value = 9900
$ws.Range("H45").Value = 2000
$ws.Range("H48").Value = 5900
$ws.Range("H49").Value = 34900
$ws.Range("H52").Value = 6700
